$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: replace the (incorrect) date with an approximate text label, right-aligned
# using the existing mm/yyyy number format.
$ws.Range("A2").Value = "~06/2017"
$ws.Range("A2").NumberFormat = "mm/yyyy"
$ws.Range("A2").HorizontalAlignment = -4152

# Extend the running "days since last" formulas down through row 28
$ws.Range("B27").Formula = "=C27-C26"
$ws.Range("B28").Formula = "=C28-C27"

# Correct the page count recorded for 2019-07 and add the new 2019-08 entry
$ws.Range("C28").Value = 582
$ws.Range("B29").Formula = "=C29-C28"
$ws.Range("C29").Value = 606

# Move the active selection like the author left it
$ws.Range("D32").Select() | Out-Null
